# Apply production_data.xlsx changes:
#  - Rows 61, 68, 84: Status (K) "Paused" -> "Finished", and stamp an End Timestamp (T)
#  - Append new task rows 102-110 at the bottom of the log
#  - Sheet dimension grows from A1:T101 to A1:T110 automatically as cells are filled

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force a numeric-looking value to be stored as text (matches the
# original export's inline-string typing for columns like H/I/M/etc.)
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row 61: finish the paused carpentry task and stamp its end time ---
$ws.Range("K61").Value = "Finished"
$ws.Range("T61").Value = "2024-07-15 00:20"

# --- Row 68: finish the paused flooring task and stamp its end time ---
$ws.Range("K68").Value = "Finished"
$ws.Range("T68").Value = "2024-07-15 00:20"

# --- Row 84: finish the paused cleanup task and stamp its end time ---
$ws.Range("K84").Value = "Finished"
$ws.Range("T84").Value = "2024-07-15 00:57"

# --- New row 102 ---
$ws.Range("A102").Value = "ec7bec93-68de-45e3-a277-503428236ec5"
$ws.Range("B102").Value = "2024-07-15 00:20"
$ws.Range("C102").Value = 32
$ws.Range("D102").Value = "CESAR VILLARROEL"
$ws.Range("E102").Value = "Carlos Astorga"
$ws.Range("F102").Value = "CARPINTERIA"
$ws.Range("G102").Value = "Puyaral"
Set-TextValue $ws.Range("H102") "25"
Set-TextValue $ws.Range("I102") "1"
$ws.Range("J102").Value = "Guardapolvos y pilastras"
$ws.Range("K102").Value = "Finished"
$ws.Range("L102").Value = 1
$ws.Range("M102").Value = "L1"
$ws.Range("T102").Value = "2024-07-15 00:20"

# --- New row 103 ---
$ws.Range("A103").Value = "06b0cafd-b600-4510-a1af-156aadf6f614"
$ws.Range("B103").Value = "2024-07-15 00:20"
$ws.Range("C103").Value = 32
$ws.Range("D103").Value = "CESAR VILLARROEL"
$ws.Range("E103").Value = "Carlos Astorga"
$ws.Range("F103").Value = "CARPINTERIA"
$ws.Range("G103").Value = "Puyaral"
Set-TextValue $ws.Range("H103") "25"
Set-TextValue $ws.Range("I103") "1"
$ws.Range("J103").Value = "Instalación OSB"
$ws.Range("K103").Value = "Finished"
$ws.Range("L103").Value = 1
$ws.Range("M103").Value = "L1"
$ws.Range("T103").Value = "2024-07-15 00:20"

# --- New row 104 ---
$ws.Range("A104").Value = "8892e157-d06e-4ad5-94f7-288555094449"
$ws.Range("B104").Value = "2024-07-15 00:20"
$ws.Range("C104").Value = 32
$ws.Range("D104").Value = "CESAR VILLARROEL"
$ws.Range("E104").Value = "Carlos Astorga"
$ws.Range("F104").Value = "CARPINTERIA"
$ws.Range("G104").Value = "Puyaral"
Set-TextValue $ws.Range("H104") "25"
Set-TextValue $ws.Range("I104") "1"
$ws.Range("J104").Value = "Instalación OSB"
$ws.Range("K104").Value = "Finished"
$ws.Range("L104").Value = 1
$ws.Range("M104").Value = "L1"
$ws.Range("T104").Value = "2024-07-15 00:20"

# --- New row 105 ---
$ws.Range("A105").Value = "a2ef94a5-c02a-4c49-b85c-7787d9182cd6"
$ws.Range("B105").Value = "2024-07-15 00:21"
$ws.Range("C105").Value = 32
$ws.Range("D105").Value = "CESAR VILLARROEL"
$ws.Range("E105").Value = "Carlos Astorga"
$ws.Range("F105").Value = "CARPINTERIA"
$ws.Range("G105").Value = "Puyaral"
Set-TextValue $ws.Range("H105") "25"
Set-TextValue $ws.Range("I105") "1"
$ws.Range("J105").Value = "Instalación OSB"
$ws.Range("K105").Value = "Finished"
$ws.Range("L105").Value = 1
$ws.Range("M105").Value = "L1"
$ws.Range("T105").Value = "2024-07-15 00:21"

# --- New row 106 ---
$ws.Range("A106").Value = "b0e9d0ae-e261-4ba4-be9f-813fa269b1a6"
$ws.Range("B106").Value = "2024-07-15 00:21"
$ws.Range("C106").Value = 32
$ws.Range("D106").Value = "CESAR VILLARROEL"
$ws.Range("E106").Value = "Carlos Astorga"
$ws.Range("F106").Value = "CARPINTERIA"
$ws.Range("G106").Value = "Puyaral"
Set-TextValue $ws.Range("H106") "25"
Set-TextValue $ws.Range("I106") "1"
$ws.Range("J106").Value = "Instalación OSB"
$ws.Range("K106").Value = "en proceso"
$ws.Range("L106").Value = 1
$ws.Range("M106").Value = "L1"

# --- New row 107 ---
$ws.Range("A107").Value = "ac8b686b-652e-41cd-912a-be7752f2f374"
$ws.Range("B107").Value = "2024-07-15 00:58"
$ws.Range("C107").Value = 12
$ws.Range("D107").Value = "FRANCISCO DIAZ"
$ws.Range("E107").Value = "Flor Sanhueza"
$ws.Range("F107").Value = "SELLOS Y ASEO"
$ws.Range("G107").Value = "Puyaral"
Set-TextValue $ws.Range("H107") "25"
Set-TextValue $ws.Range("I107") "1"
$ws.Range("J107").Value = "Aseo entrega"
$ws.Range("K107").Value = "Paused"
$ws.Range("L107").Value = 1
$ws.Range("M107").Value = "L1"
$ws.Range("N107").Value = "2024-07-15 00:59"
$ws.Range("O107").Value = "Final del día"

# --- New row 108 ---
$ws.Range("A108").Value = "3c1e37fd-6866-4c4e-b4df-fdf3aa9e9c1c"
$ws.Range("B108").Value = "2024-07-15 09:08"
$ws.Range("C108").Value = 25
$ws.Range("D108").Value = "CESAR VILLARROEL"
$ws.Range("E108").Value = "CELSO MARTINEZ"
$ws.Range("F108").Value = "CARPINTERIA"
$ws.Range("G108").Value = "Las Bandurrias"
Set-TextValue $ws.Range("H108") "15"
Set-TextValue $ws.Range("I108") "1"
$ws.Range("J108").Value = "Instalación de puerta interior"
$ws.Range("K108").Value = "Finished"
$ws.Range("L108").Value = 1
$ws.Range("M108").Value = "L1"
$ws.Range("T108").Value = "2024-07-15 09:08"

# --- New row 109 ---
$ws.Range("A109").Value = "240f279e-469b-45fd-b098-cdcbc15480b9"
$ws.Range("B109").Value = "2024-07-15 09:08"
$ws.Range("C109").Value = 25
$ws.Range("D109").Value = "CESAR VILLARROEL"
$ws.Range("E109").Value = "Celso Martinez"
$ws.Range("F109").Value = "CARPINTERIA"
$ws.Range("G109").Value = "Las Bandurrias"
Set-TextValue $ws.Range("H109") "15"
Set-TextValue $ws.Range("I109") "1"
$ws.Range("J109").Value = "Revestimiento Siding (planchas)"
$ws.Range("K109").Value = "en proceso"
$ws.Range("L109").Value = 1
$ws.Range("M109").Value = "L1"

# --- New row 110 ---
$ws.Range("A110").Value = "a2a02727-8a93-43a4-a62c-8f836eb05283"
$ws.Range("B110").Value = "2024-07-15 09:09"
$ws.Range("C110").Value = 12
$ws.Range("D110").Value = "FRANCISCO DIAZ"
$ws.Range("E110").Value = "FLOR SANHUEZA"
$ws.Range("F110").Value = "SELLOS Y ASEO"
$ws.Range("G110").Value = "Bosquemar"
Set-TextValue $ws.Range("H110") "15"
Set-TextValue $ws.Range("I110") "4"
$ws.Range("J110").Value = "Aseo entrega"
$ws.Range("K110").Value = "Finished"
$ws.Range("L110").Value = 4
$ws.Range("M110").Value = "L2"
$ws.Range("T110").Value = "2024-07-15 09:09"
